$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F4").Value = 4
$ws.Range("F7").Value = 10
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -4
